# Refresh the best-month-per-route dataset with the test-1 results (see commit message:
# "Finish test 1; improve code based on test 1 results"). Route order was reshuffled,
# several routes seasonal stats/dates were recomputed, and a third "California Plus"
# sample (day 312) was appended as new row 34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    [pscustomobject]@{ Row=2; Route="East Canada"; Day=177; Temperate=0.84; Mean15d=0.82; Half=0; Best=$true; Date="06-26" }
    [pscustomobject]@{ Row=3; Route="East Canada"; Day=257; Temperate=0.86; Mean15d=0.86; Half=1; Best=$true; Date="09-14" }
    [pscustomobject]@{ Row=4; Route="Wyoming Plus"; Day=157; Temperate=0.71; Mean15d=0.6899999999999999; Half=0; Best=$true; Date="06-06" }
    [pscustomobject]@{ Row=5; Route="Wyoming Plus"; Day=272; Temperate=0.67; Mean15d=0.67; Half=1; Best=$true; Date="09-29" }
    [pscustomobject]@{ Row=6; Route="Louisiana Plus"; Day=97; Temperate=0.78; Mean15d=0.76; Half=0; Best=$true; Date="04-07" }
    [pscustomobject]@{ Row=7; Route="Louisiana Plus"; Day=337; Temperate=0.72; Mean15d=0.71; Half=1; Best=$true; Date="12-03" }
    [pscustomobject]@{ Row=8; Route="Alaska State"; Day=212; Temperate=0.84; Mean15d=0.85; Half=0; Best=$true; Date="07-31" }
    [pscustomobject]@{ Row=9; Route="Alaska State"; Day=217; Temperate=0.86; Mean15d=0.85; Half=1; Best=$true; Date="08-05" }
    [pscustomobject]@{ Row=10; Route="New Mexico Plus"; Day=112; Temperate=0.66; Mean15d=0.65; Half=0; Best=$true; Date="04-22" }
    [pscustomobject]@{ Row=11; Route="New Mexico Plus"; Day=297; Temperate=0.63; Mean15d=0.64; Half=1; Best=$true; Date="10-24" }
    [pscustomobject]@{ Row=12; Route="Minnesota Plus"; Day=177; Temperate=0.85; Mean15d=0.78; Half=0; Best=$true; Date="06-26" }
    [pscustomobject]@{ Row=13; Route="Minnesota Plus"; Day=252; Temperate=0.82; Mean15d=0.82; Half=1; Best=$true; Date="09-09" }
    [pscustomobject]@{ Row=14; Route="Indiana Plus"; Day=157; Temperate=0.79; Mean15d=0.76; Half=0; Best=$true; Date="06-06" }
    [pscustomobject]@{ Row=15; Route="Indiana Plus"; Day=257; Temperate=0.84; Mean15d=0.79; Half=1; Best=$true; Date="09-14" }
    [pscustomobject]@{ Row=16; Route="Northeast Plus"; Day=157; Temperate=0.79; Mean15d=0.76; Half=0; Best=$true; Date="06-06" }
    [pscustomobject]@{ Row=17; Route="Northeast Plus"; Day=277; Temperate=0.79; Mean15d=0.8; Half=1; Best=$true; Date="10-04" }
    [pscustomobject]@{ Row=18; Route="Florida State"; Day=37; Temperate=0.75; Mean15d=0.73; Half=0; Best=$true; Date="02-06" }
    [pscustomobject]@{ Row=19; Route="Florida State"; Day=362; Temperate=0.75; Mean15d=0.72; Half=1; Best=$true; Date="12-28" }
    [pscustomobject]@{ Row=20; Route="Missouri Plus"; Day=147; Temperate=0.76; Mean15d=0.76; Half=0; Best=$true; Date="05-27" }
    [pscustomobject]@{ Row=21; Route="Missouri Plus"; Day=282; Temperate=0.78; Mean15d=0.74; Half=1; Best=$true; Date="10-09" }
    [pscustomobject]@{ Row=22; Route="Georgia Plus"; Day=117; Temperate=0.83; Mean15d=0.79; Half=0; Best=$true; Date="04-27" }
    [pscustomobject]@{ Row=23; Route="Georgia Plus"; Day=302; Temperate=0.79; Mean15d=0.79; Half=1; Best=$true; Date="10-29" }
    [pscustomobject]@{ Row=24; Route="Hawaii State"; Day=72; Temperate=0.8100000000000001; Mean15d=0.83; Half=0; Best=$true; Date="03-13" }
    [pscustomobject]@{ Row=25; Route="Hawaii State"; Day=365; Temperate=0.67; Mean15d=0.67; Half=1; Best=$true; Date="12-31" }
    [pscustomobject]@{ Row=26; Route="Near DC"; Day=137; Temperate=0.77; Mean15d=0.78; Half=0; Best=$true; Date="05-17" }
    [pscustomobject]@{ Row=27; Route="Near DC"; Day=297; Temperate=0.8; Mean15d=0.82; Half=1; Best=$true; Date="10-24" }
    [pscustomobject]@{ Row=28; Route="Oregon Plus"; Day=177; Temperate=0.78; Mean15d=0.8100000000000001; Half=0; Best=$true; Date="06-26" }
    [pscustomobject]@{ Row=29; Route="Oregon Plus"; Day=267; Temperate=0.82; Mean15d=0.85; Half=1; Best=$true; Date="09-24" }
    [pscustomobject]@{ Row=30; Route="Puerto Rico"; Day=32; Temperate=0.14; Mean15d=0.2; Half=0; Best=$true; Date="02-01" }
    [pscustomobject]@{ Row=31; Route="Puerto Rico"; Day=365; Temperate=0.13; Mean15d=0.1; Half=1; Best=$true; Date="12-31" }
    [pscustomobject]@{ Row=32; Route="California Plus"; Day=147; Temperate=0.77; Mean15d=0.79; Half=0; Best=$true; Date="05-27" }
    [pscustomobject]@{ Row=33; Route="California Plus"; Day=311; Temperate=0.76; Mean15d=0.76; Half=1; Best=$true; Date="11-07" }
    [pscustomobject]@{ Row=34; Route="California Plus"; Day=312; Temperate=0.76; Mean15d=0.76; Half=1; Best=$true; Date="11-08" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Route
    $ws.Cells.Item($r.Row, 2).Value = $r.Day
    $ws.Cells.Item($r.Row, 3).Value = $r.Temperate
    $ws.Cells.Item($r.Row, 4).Value = $r.Mean15d
    $ws.Cells.Item($r.Row, 5).Value = $r.Half
    $ws.Cells.Item($r.Row, 6).Value = $r.Best
    $ws.Cells.Item($r.Row, 7).Value = $r.Date
}

Write-Host "Updated $($rows.Count) data rows (A2:G34) on sheet $($ws.Name)."
